$wb = $excel.ActiveWorkbook

# --- Rename the "Include from RoleClass" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from RoleClass")
$wsInclude.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row before the "Description" row (currently row 11),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Copy formatting from the row below (the shifted "Description" row) so the new
# row matches the existing data-row style instead of getting a blank default style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
